$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 was a blank placeholder row (only H8 had the "❌" marker). Replace
# it in place with the "566 - Reshape the matrix" entry; rows below are
# left untouched.

# Pull the per-column formatting (fills / fonts / wrap text) from row 7,
# which already carries the look the new data row should have.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null

$ws.Range("D7").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null

$ws.Range("H7").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null

$ws.Range("I7").Copy() | Out-Null
$ws.Range("I8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$ws.Range("A8").Value = 566
$ws.Range("B8").Value = "Reshape the matrix"
$ws.Range("C8").Value = "Array"
$ws.Range("D8").Value = "Array, Matrix, Simulation"
$ws.Range("F8").Value = "Easy"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "✅"
$ws.Range("I8").Value = "Given O(n^2) sol and didn't see any other sol"
$ws.Rows.Item(8).RowHeight = 30

# Restore the active selection to E13 like the recorded session.
$ws.Range("E13").Select() | Out-Null
